# Applies the coinranking.com price-table refresh described in the commit
# "Updated cryptos list on Tue Jul 11 03:31:38 UTC 2023 with GitHub Actions".
# Column D (Price) holds numeric-looking text (e.g. "1.001", "30.431.87")
# that must stay text, exactly as scraped - so for those cells we briefly
# force a Text number format before writing the value, then restore the
# cell to the Normal style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.431.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.76%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.877.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.72%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.84%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4756'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.93%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2892'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.31%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06506'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.84%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.13%  '

# Row 12
$ws.Range("E12").Value = '  +8.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '96.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.879.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.29%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.112'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.444.76'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.88%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007541'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '

# Row 20
$ws.Range("E20").Value = '  +0.00%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.119.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.261'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.157'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.94%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.262'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.87%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.953'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.31%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.373'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09971'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.41%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.518'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.312'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.056'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.56%  '

# Row 34
$ws.Range("E34").Value = '  +2.09%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.123'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01865'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.737'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.342'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.18%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.00%  '

# Row 42
$ws.Range("E42").Value = '  +1.89%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4162'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.54%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8367'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.22%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.314'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.51%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.069'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.53%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.74%  '

# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '925.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.18%  '

# Row 51
$ws.Range("E51").Value = '  +0.67%  '
